# The commit swaps the content of the two theme parts in the package:
#   ppt/theme/theme1.xml (used by the notes master)  <->  Integral / "Red Violet" theme
#   ppt/theme/theme2.xml (used by the slide master / presentation) <->  default "Office Theme"
#
# i.e. the deck's visible design switches from the pink/purple "Integral" theme
# to the plain "Office Theme" palette, while the notes master keeps the Integral
# colours. The PowerPoint object model only exposes the *active* design's theme
# (the one driving the slide master / slides), reachable as
# $p.SlideMaster.Theme.ThemeColorScheme - so we repoint those 12 colour slots to
# the "Office Theme" values that the diff puts in theme2.xml.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$tcs = $m.Theme.ThemeColorScheme

# Office Theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
# PowerPoint's ColorFormat.RGB is a BGR-packed integer (0xBBGGRR), i.e. the
# same value VBA's RGB(r,g,b) = r + g*256 + b*65536 would produce.
$tcs.Item(1).RGB  = 0x000000    # dk1      000000
$tcs.Item(2).RGB  = 0xFFFFFF    # lt1      FFFFFF
$tcs.Item(3).RGB  = 0x6A5444    # dk2      44546A
$tcs.Item(4).RGB  = 0xE6E6E7    # lt2      E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B    # accent1  5B9BD5
$tcs.Item(6).RGB  = 0x317DED    # accent2  ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5    # accent3  A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF    # accent4  FFC000
$tcs.Item(9).RGB  = 0xC47244    # accent5  4472C4
$tcs.Item(10).RGB = 0x47AD70    # accent6  70AD47
$tcs.Item(11).RGB = 0xC16305    # hlink    0563C1
$tcs.Item(12).RGB = 0x724F95    # folHlink 954F72
